$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "39.806.14"
Set-TextValue "E2" "  +1.03%  "
Set-TextValue "D3" "2.220.32"
Set-TextValue "E3" "  +0.51%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "292.56"
Set-TextValue "E5" "  -1.78%  "
Set-TextValue "D6" "86.23"
Set-TextValue "E6" "  +4.61%  "
Set-TextValue "D7" "0.515"
Set-TextValue "E7" "  +0.77%  "
Set-TextValue "D8" "1.00"
Set-TextValue "E8" "  -0.08%  "
Set-TextValue "D9" "0.473"
Set-TextValue "E9" "  +0.84%  "
Set-TextValue "D10" "30.84"
Set-TextValue "E10" "  +5.92%  "
Set-TextValue "D11" "0.0786"
Set-TextValue "E11" "  +1.53%  "
Set-TextValue "D12" "47.13"
Set-TextValue "E12" "  -0.82%  "
Set-TextValue "E13" "  +1.59%  "
Set-TextValue "D14" "6.36"
Set-TextValue "E14" "  +1.11%  "
Set-TextValue "D15" "2.559.89"
Set-TextValue "E15" "  +0.21%  "
Set-TextValue "D16" "14.11"
Set-TextValue "E16" "  +0.27%  "
Set-TextValue "D17" "2.220.02"
Set-TextValue "E17" "  +0.85%  "
Set-TextValue "D18" "0.729"
Set-TextValue "E18" "  +2.23%  "
Set-TextValue "D19" "39.763.42"
Set-TextValue "E19" "  +1.19%  "
Set-TextValue "D20" "0.0₃0883"
Set-TextValue "E20" "  +0.80%  "
Set-TextValue "D21" "11.09"
Set-TextValue "E21" "  +7.49%  "
Set-TextValue "D22" "5.81"
Set-TextValue "E22" "  +1.65%  "
Set-TextValue "D23" "65.60"
Set-TextValue "E23" "  +0.85%  "
Set-TextValue "D24" "236.04"
Set-TextValue "E24" "  +4.20%  "
Set-TextValue "D25" "1.00"
Set-TextValue "E25" "  +0.15%  "
Set-TextValue "D26" "2.47"
Set-TextValue "E26" "  +2.49%  "
Set-TextValue "E27" "  +2.78%  "
Set-TextValue "D28" "22.75"
Set-TextValue "E28" "  +0.57%  "
Set-TextValue "D29" "2.20"
Set-TextValue "E29" "  +1.72%  "
Set-TextValue "D30" "9.26"
Set-TextValue "E30" "  +1.77%  "
Set-TextValue "D31" "32.94"
Set-TextValue "E31" "  +3.67%  "
Set-TextValue "D32" "151.51"
Set-TextValue "E32" "  +1.85%  "
Set-TextValue "D33" "0.998"
Set-TextValue "E33" "  -0.24%  "
Set-TextValue "D34" "4.95"
Set-TextValue "E34" "  +2.56%  "
Set-TextValue "D35" "0.0719"
Set-TextValue "E35" "  +3.67%  "
Set-TextValue "E36" "  +1.72%  "
Set-TextValue "E37" "  +7.32%  "
Set-TextValue "D38" "16.06"
Set-TextValue "E38" "  +5.47%  "
Set-TextValue "D39" "0.112"
Set-TextValue "E39" "  +1.95%  "
Set-TextValue "D40" "0.0994"
Set-TextValue "E40" "  +3.19%  "
Set-TextValue "D41" "1.70"
Set-TextValue "E41" "  +3.69%  "
Set-TextValue "D42" "3.79"
Set-TextValue "E42" "  +4.48%  "
Set-TextValue "D43" "2.060.14"
Set-TextValue "E43" "  +8.55%  "
Set-TextValue "D44" "18.03"
Set-TextValue "E44" "  +11.66%  "
Set-TextValue "B45" "VeChain"
Set-TextValue "C45" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D45" "0.0268"
Set-TextValue "E45" "  +3.64%  "
Set-TextValue "B46" "ApeXProtocol"
Set-TextValue "C46" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D46" "2.10"
Set-TextValue "E46" "  +0.76%  "
Set-TextValue "D47" "9.95"
Set-TextValue "E47" "  +10.44%  "
Set-TextValue "D48" "2.62"
Set-TextValue "E48" "  +0.32%  "
Set-TextValue "D49" "2.432.54"
Set-TextValue "E49" "  +0.38%  "
Set-TextValue "D50" "71.73"
Set-TextValue "E50" "  +1.43%  "
Set-TextValue "D51" "89.17"
Set-TextValue "E51" "  +2.35%  "
